# "change subject to subject_topic"
#
# The HP image metadata template has a header row (row 1) of field names.
# The column that held the literal header text "subject" is renamed to
# "subject_topic" (column C in this sheet). Shared-string bookkeeping
# (removing the now-unused "subject" entry / appending "subject_topic")
# is handled automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C1").Value = "subject_topic"

# Cosmetic touches from the same save (harmless / match author's final view
# state): the header row was shortened from its old wrapped height and the
# active selection ended up on V5 after the edit.
$ws.Rows("1:1").RowHeight = 17
$ws.Range("V5").Select()
